$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column L: copy the 2020 column (K) into L
$ws.Range("L3").Value = 2020
$ws.Range("L4").Value = 6.18

# Update the style of the new cells to match column K (style copied via Range.Style mapping
# handled below by copying formatting from K3/K4)
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)  # xlPasteFormats

# Update the selection to match the diff
$ws.Range("L10").Select()
